# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" sheet (fund holdings detail, same shape as the
# other quarterly sheets) immediately before the "总计" (totals) sheet, and
# adds a corresponding summary row at the top of "总计".

function Set-HeaderCell($cell, $text) {
    $cell.Value = $text
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
}

function Set-IndexCell($cell, $n) {
    # Mimics the bold/centered "row index" style used in column A of every
    # quarterly sheet (0-based position within that sheet).
    $cell.Value = $n
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
}

function Set-TextValueCell($cell, $text) {
    # Forces a numeric-looking string (e.g. "0.18") to be stored as text
    # instead of being auto-coerced to a number, then strips the
    # number-format override back off so no extra style sticks around.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right before "总计".
# ---------------------------------------------------------------------
$totalSheetRef = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheetRef, $null)
$q1.Name = "2022-Q1"

# The "Before" reference captured above can go stale once the sheet
# collection is mutated - re-resolve "总计" by name now that the insert
# (and rename) are done, so later writes land on the right sheet object.
$totalSheet = $wb.Worksheets.Item("总计")

# Header row (B1:H1) - same columns as the other quarterly detail sheets.
$h1 = $q1.Cells.Item(1, 2)
Set-HeaderCell $h1 "基金代码"
$h2 = $q1.Cells.Item(1, 3)
Set-HeaderCell $h2 "基金名称"
$h3 = $q1.Cells.Item(1, 4)
Set-HeaderCell $h3 "基金规模"
$h4 = $q1.Cells.Item(1, 5)
Set-HeaderCell $h4 "股票总仓位"
$h5 = $q1.Cells.Item(1, 6)
Set-HeaderCell $h5 "仓位占比"
$h6 = $q1.Cells.Item(1, 7)
Set-HeaderCell $h6 "持有市值(亿元)"
$h7 = $q1.Cells.Item(1, 8)
Set-HeaderCell $h7 "仓位排名"

# Data rows: code, name, scale, total stock position, position %,
# held market value (yi), position rank.
$q1Codes  = @("001303", "002323")
$q1Names  = @("银华稳利灵活配置混合A", "银华稳利灵活配置混合C")
$q1Scale  = @("0.18", "0.12")
$q1Pos    = @("28.88", "28.88")
$q1Pct    = @("0.61", "0.61")
$q1Value  = @("0.0011", "0.0007")
$q1Rank   = @(8, 8)

for ($i = 0; $i -lt 2; $i++) {
    $r = $i + 2

    $idxCell = $q1.Cells.Item($r, 1)
    Set-IndexCell $idxCell $i

    $codeCell = $q1.Cells.Item($r, 2)
    $codeText = $q1Codes[$i]
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $codeText
    $codeCell.ClearFormats()

    $nameCell = $q1.Cells.Item($r, 3)
    $nameCell.Value = $q1Names[$i]

    $scaleCell = $q1.Cells.Item($r, 4)
    $scaleText = $q1Scale[$i]
    Set-TextValueCell $scaleCell $scaleText

    $posCell = $q1.Cells.Item($r, 5)
    $posText = $q1Pos[$i]
    Set-TextValueCell $posCell $posText

    $pctCell = $q1.Cells.Item($r, 6)
    $pctText = $q1Pct[$i]
    Set-TextValueCell $pctCell $pctText

    $valueCell = $q1.Cells.Item($r, 7)
    $valueText = $q1Value[$i]
    Set-TextValueCell $valueCell $valueText

    $rankCell = $q1.Cells.Item($r, 8)
    $rankCell.Value = $q1Rank[$i]
}

# ---------------------------------------------------------------------
# 2. Update "总计": insert a 2022-Q1 summary row at the top, shifting the
#    existing quarters down by one and renumbering the index column.
# ---------------------------------------------------------------------
$totalLabels = @("2022-Q1", "2021-Q4", "2021-Q3", "2021-Q2", "2021-Q1", "2020-Q4")
$totalCounts = @(2, 14, 18, 10, 14, 8)
$totalValues = @(0, 2.21, 5.66, 3.64, 5.55, 0.21)

for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2

    $idxCell = $totalSheet.Cells.Item($r, 1)
    $idxCell.Value = $i

    $labelCell = $totalSheet.Cells.Item($r, 2)
    $labelCell.Value = $totalLabels[$i]

    $countCell = $totalSheet.Cells.Item($r, 3)
    $countCell.Value = $totalCounts[$i]

    $valueCell = $totalSheet.Cells.Item($r, 4)
    $valueCell.Value = $totalValues[$i]
}

# The "总计" sheet grew by one row (row 7 is brand new), so it has no
# pre-existing style to inherit - apply the same bold/centered index style
# used by the rest of column A.
$newIdxCell = $totalSheet.Cells.Item(7, 1)
$newIdxCell.Font.Bold = $true
$newIdxCell.HorizontalAlignment = -4108
$newIdxCell.VerticalAlignment = -4160

# Adding a sheet makes it the active one; restore the original active tab
# ("2020-Q4" is the first sheet) since the edit itself doesn't change which
# sheet was selected.
$firstSheet = $wb.Worksheets.Item("2020-Q4")
$firstSheet.Activate()
